# Auto-generated update of cryptos.xlsx Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Force the value to be stored as text (matching the original inlineStr cell type),
    # even when it looks like a plain number (e.g. "313.83" or "1.000").
    if ($value -match '^-?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}

Set-TextCell $ws.Cells.Item(2, 4) "27.317.51"
$ws.Cells.Item(2, 5).Value = "  +0.98%  "

Set-TextCell $ws.Cells.Item(3, 4) "1.776.59"
$ws.Cells.Item(3, 5).Value = "  +3.98%  "

Set-TextCell $ws.Cells.Item(4, 4) "1.000"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "

Set-TextCell $ws.Cells.Item(5, 4) "313.83"
$ws.Cells.Item(5, 5).Value = "  +2.06%  "

Set-TextCell $ws.Cells.Item(6, 4) "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.04%  "

Set-TextCell $ws.Cells.Item(7, 4) "0.5250"
$ws.Cells.Item(7, 5).Value = "  +11.16%  "

Set-TextCell $ws.Cells.Item(8, 4) "0.3612"
$ws.Cells.Item(8, 5).Value = "  +5.54%  "

Set-TextCell $ws.Cells.Item(9, 4) "42.50"
$ws.Cells.Item(9, 5).Value = "  +1.25%  "

Set-TextCell $ws.Cells.Item(10, 4) "0.07372"
$ws.Cells.Item(10, 5).Value = "  +1.53%  "

$ws.Cells.Item(11, 5).Value = "  +5.60%  "

Set-TextCell $ws.Cells.Item(12, 4) "1.002"
$ws.Cells.Item(12, 5).Value = "  +0.00%  "

Set-TextCell $ws.Cells.Item(13, 4) "20.53"
$ws.Cells.Item(13, 5).Value = "  +3.98%  "

Set-TextCell $ws.Cells.Item(14, 4) "6.068"
$ws.Cells.Item(14, 5).Value = "  +4.21%  "

Set-TextCell $ws.Cells.Item(15, 4) "1.775.89"
$ws.Cells.Item(15, 5).Value = "  +4.10%  "

Set-TextCell $ws.Cells.Item(16, 4) "6.979"
$ws.Cells.Item(16, 5).Value = "  +2.42%  "

Set-TextCell $ws.Cells.Item(17, 4) "88.49"
$ws.Cells.Item(17, 5).Value = "  -0.42%  "

Set-TextCell $ws.Cells.Item(18, 4) "0.00001044"
$ws.Cells.Item(18, 5).Value = "  +1.05%  "

Set-TextCell $ws.Cells.Item(19, 4) "0.06421"
$ws.Cells.Item(19, 5).Value = "  +1.10%  "

Set-TextCell $ws.Cells.Item(20, 4) "1.001"
$ws.Cells.Item(20, 5).Value = "  -0.04%  "

$ws.Cells.Item(21, 5).Value = "  +1.96%  "

Set-TextCell $ws.Cells.Item(22, 4) "5.843"
$ws.Cells.Item(22, 5).Value = "  +4.53%  "

Set-TextCell $ws.Cells.Item(23, 4) "27.380.69"
$ws.Cells.Item(23, 5).Value = "  +1.13%  "

Set-TextCell $ws.Cells.Item(24, 4) "11.33"
$ws.Cells.Item(24, 5).Value = "  +4.72%  "

Set-TextCell $ws.Cells.Item(25, 4) "2.074"
$ws.Cells.Item(25, 5).Value = "  -1.18%  "

Set-TextCell $ws.Cells.Item(26, 4) "153.80"
$ws.Cells.Item(26, 5).Value = "  -1.21%  "

Set-TextCell $ws.Cells.Item(27, 4) "20.09"
$ws.Cells.Item(27, 5).Value = "  +2.91%  "

Set-TextCell $ws.Cells.Item(28, 4) "2.354"
$ws.Cells.Item(28, 5).Value = "  +13.70%  "

Set-TextCell $ws.Cells.Item(29, 4) "1.980.79"
$ws.Cells.Item(29, 5).Value = "  +4.19%  "

Set-TextCell $ws.Cells.Item(30, 4) "121.38"
$ws.Cells.Item(30, 5).Value = "  +1.92%  "

Set-TextCell $ws.Cells.Item(31, 4) "1.065"
$ws.Cells.Item(31, 5).Value = "  +5.91%  "

Set-TextCell $ws.Cells.Item(32, 4) "0.09791"
$ws.Cells.Item(32, 5).Value = "  +7.23%  "

Set-TextCell $ws.Cells.Item(33, 4) "5.555"

Set-TextCell $ws.Cells.Item(34, 4) "3.619"
$ws.Cells.Item(34, 5).Value = "  +1.12%  "

Set-TextCell $ws.Cells.Item(35, 4) "0.02233"
$ws.Cells.Item(35, 5).Value = "  +2.43%  "

Set-TextCell $ws.Cells.Item(36, 4) "0.05969"
$ws.Cells.Item(36, 5).Value = "  +2.66%  "

Set-TextCell $ws.Cells.Item(37, 4) "11.21"
$ws.Cells.Item(37, 5).Value = "  +1.82%  "

Set-TextCell $ws.Cells.Item(38, 4) "4.856"
$ws.Cells.Item(38, 5).Value = "  +3.04%  "

Set-TextCell $ws.Cells.Item(39, 4) "0.2026"
$ws.Cells.Item(39, 5).Value = "  +2.42%  "

Set-TextCell $ws.Cells.Item(40, 4) "0.6140"
$ws.Cells.Item(40, 5).Value = "  +4.96%  "

Set-TextCell $ws.Cells.Item(41, 4) "1.428"
$ws.Cells.Item(41, 5).Value = "  +2.67%  "

Set-TextCell $ws.Cells.Item(42, 4) "8.102"
$ws.Cells.Item(42, 5).Value = "  +8.92%  "

Set-TextCell $ws.Cells.Item(43, 4) "1.145"
$ws.Cells.Item(43, 5).Value = "  +4.32%  "

$ws.Cells.Item(44, 5).Value = "  +5.00%  "

Set-TextCell $ws.Cells.Item(45, 4) "0.5771"
$ws.Cells.Item(45, 5).Value = "  +2.77%  "

Set-TextCell $ws.Cells.Item(46, 4) "3.628"
$ws.Cells.Item(46, 5).Value = "  +2.07%  "

Set-TextCell $ws.Cells.Item(47, 4) "121.49"
$ws.Cells.Item(47, 5).Value = "  +3.44%  "

Set-TextCell $ws.Cells.Item(48, 4) "1.887"
$ws.Cells.Item(48, 5).Value = "  +3.05%  "

Set-TextCell $ws.Cells.Item(49, 4) "1.110"
$ws.Cells.Item(49, 5).Value = "  +2.79%  "

Set-TextCell $ws.Cells.Item(50, 4) "0.06711"
$ws.Cells.Item(50, 5).Value = "  +1.31%  "

Set-TextCell $ws.Cells.Item(51, 4) "70.88"
$ws.Cells.Item(51, 5).Value = "  +2.00%  "
